$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-224) holds a "Förändrad" (Changed) date that was
# bumped by one day, from serial 45188 (2023-09-19) to serial 45189
# (2023-09-20), for every data row in the sheet.
$ws.Range("C2:C224").Value = 45189
